$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.905.75"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "2.529.74"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'590.51"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").Value = "'172.88"
$ws.Range("E6").Value = "  -1.59%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.523"
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("D9").Value = "2.529.92"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("E10").Value = "  -2.67%  "
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("D12").Value = "'0.343"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "'5.01"
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("D14").Value = "'26.32"
$ws.Range("E14").Value = "  -1.73%  "
$ws.Range("D15").Value = "2.989.89"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("D17").Value = "68.071.66"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "2.536.73"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").Value = "'11.89"
$ws.Range("E19").Value = "  +3.97%  "
$ws.Range("D20").Value = "'8.02"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "'364.73"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").Value = "'4.14"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").Value = "'1.32"
$ws.Range("E23").Value = "  +32.71%  "
$ws.Range("D24").Value = "'4.51"
$ws.Range("E24").Value = "  -2.88%  "
$ws.Range("D25").Value = "'71.96"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("E27").Value = "  -4.84%  "
$ws.Range("D28").Value = "'9.84"
$ws.Range("E28").Value = "  -4.39%  "
$ws.Range("D29").Value = "2.664.95"
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").Value = "0.0₃0951"
$ws.Range("E30").Value = "  -3.59%  "
$ws.Range("D31").Value = "'530.58"
$ws.Range("E31").Value = "  -3.84%  "
$ws.Range("D32").Value = "'8.24"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").Value = "'1.29"
$ws.Range("E34").Value = "  -4.03%  "
$ws.Range("D35").Value = "'0.128"
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'159.77"
$ws.Range("E37").Value = "  +2.51%  "
$ws.Range("D38").Value = "'19.29"
$ws.Range("E38").Value = "  +2.87%  "
$ws.Range("D39").Value = "'1.43"
$ws.Range("D40").Value = "'18.60"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "'1.77"
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("D42").Value = "'5.08"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("D43").Value = "'0.345"
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "'2.45"
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("D46").Value = "'39.32"
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("D47").Value = "'148.01"
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("D48").Value = "'0.552"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("D49").Value = "'3.69"
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("D50").Value = "0.0₆0274"
$ws.Range("E50").Value = "  -2.34%  "
$ws.Range("D51").Value = "'1.70"
$ws.Range("E51").Value = "  +0.89%  "
